$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "MMU & AMS" 3mf entry (row 24) is being dropped from the parts list.
# Select the row first (mirrors the interactive workflow of right-click ->
# Delete on the row header) and then delete it outright, which shifts every
# row below it up by one.
$ws.Range("A24:XFD24").Select()
$ws.Rows("24:24").Delete()

# Re-assign the formulas for the surviving data range as whole-range writes
# (rather than cell-by-cell) so the workbook re-collapses them back into
# shared formulas spanning the new A4:K35 extent, matching how Excel stores
# a fill-down/fill-right formula block.
$ws.Range("K2:K35").Formula = '="|"&$D2&"|"&$E2&"|"&$F2&"|"&$G2&"|"&$H2&"|"&$I2&"|"&$J2'
$ws.Range("B4:B35").Formula = '=_xlfn.TEXTAFTER($A4,"/",-1)'
$ws.Range("C4:C35").Formula = '=RIGHT(_xlfn.TEXTBEFORE($A4,"/",-1),LEN(_xlfn.TEXTBEFORE($A4,"/",-1))-7)'
$ws.Range("D4:D35").Formula = '="["&LEFT($B4,LEN($B4)-4)&"](../STLs/"&SUBSTITUTE($C4," ","%20")&"/"&SUBSTITUTE($B4," ","%20")&")"'
$ws.Range("E4:E35").Formula = '="!["&LEFT($B4,LEN($B4)-4)&"](./images/printed_parts/"&SUBSTITUTE($C4," ","%20")&"/"&SUBSTITUTE(SUBSTITUTE($B4,".stl",".jpg")," ","%20")&")"'
$ws.Range("F4:F35").Formula = '=$C4'
$ws.Range("G4:G35").Formula = '=IF(_xlfn.REGEXTEST($B4,"\[\S*[a][^a\]]*\]"),"Accent","Main")'
$ws.Range("H4:H35").Formula = '=IF(_xlfn.REGEXTEST($B4,"\[\S*[s][^s\]]*\]"),"Yes","No")'

# Row deletion does not automatically resize the worksheet AutoFilter or the
# hidden _xlnm._FilterDatabase defined name, so refresh both to the new
# A1:K35 extent.
$ws.AutoFilterMode = $false
$ws.Range("A1:K35").AutoFilter()

$filterDatabaseName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterDatabaseName.RefersTo = "=Sheet1!`$A`$1:`$K`$35"
